$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("endeca_attributes")

# Update row 2 values to reflect the corrected datatype dictionary entry
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "REVENUE_DISTRIBUTED_DATE"
$ws.Range("C2").Value = "mdex:dateTime"
$ws.Range("E2").Value = "Revenue Distributed Date"

# Move the active selection to E2
$ws.Range("E2").Select()
